$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.660.73'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '1.634.15'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = '''211.62'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('D8').Value = '''23.12'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').Value = '''0.0862'
$ws.Range('E11').Value = '  -3.17%  '
$ws.Range('D12').Value = '1.863.92'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '1.633.01'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '''4.05'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '''0.560'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '''65.16'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').Value = '27.632.55'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '''229.83'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').Value = '''7.58'
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').Value = '''10.66'
$ws.Range('E22').Value = '  +7.12%  '
$ws.Range('E23').Value = '  +1.69%  '
$ws.Range('E24').Value = '  +3.39%  '
$ws.Range('D25').Value = '''149.08'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('E32').Value = '  -0.83%  '
$ws.Range('D33').Value = '1.462.99'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = '''1.55'
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').Value = '''0.878'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''0.558'
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0167'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').Value = '''69.14'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('D44').Value = '''2.46'
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('E45').Value = '  -4.42%  '
$ws.Range('D46').Value = '''5.37'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').Value = '1.774.49'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('D49').Value = '''87.59'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = '0.0₆0106'
$ws.Range('E50').Value = '  +6.60%  '
$ws.Range('E51').Value = '  +0.68%  '
